$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.003.96'
$ws.Range('E2').Value = '  -0.49%  '

$ws.Range('D3').Value = '2.551.03'
$ws.Range('E3').Value = '  -0.31%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.29'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.70%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.45'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.95%  '

$ws.Range('E7').Value = '  -0.23%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.548'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.72%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.81'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.70%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.32%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.75'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.54%  '

$ws.Range('E13').Value = '  +6.12%  '

$ws.Range('D14').Value = '2.947.66'
$ws.Range('E14').Value = '  -0.02%  '

$ws.Range('D15').Value = '2.583.53'
$ws.Range('E15').Value = '  +0.96%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.881'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.57%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.85'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +4.35%  '

$ws.Range('D18').Value = '43.149.04'
$ws.Range('E18').Value = '  -0.15%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.67'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.51%  '

$ws.Range('D20').Value = '0.0₃0987'
$ws.Range('E20').Value = '  +0.09%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.62'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.32%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.89'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.48%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '256.05'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.94%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.97'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.57%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.09'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.33%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.97'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -6.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.26%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.15'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.78%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.92'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.28%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.19'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.55%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.04'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.27%  '

$ws.Range('E32').Value = '  +2.67%  '

$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.75'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.28%  '

$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.16'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.26%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0806'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.75%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.32'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.84'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +12.53%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.84'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +10.35%  '

$ws.Range('E39').Value = '  -1.05%  '

$ws.Range('E40').Value = '  -0.66%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.44'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.85%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.88'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.17%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0306'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.56%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.089.90'
$ws.Range('E44').Value = '  +1.14%  '

$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.04'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +28.35%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.20%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.79'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.90%  '

$ws.Range('E48').Value = '  +2.69%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '75.55'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +8.49%  '

$ws.Range('D50').Value = '2.804.61'
$ws.Range('E50').Value = '  +0.04%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '103.75'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.72%  '
